$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the data row that holds the daily spot price record.
# Update the date (A2) and all hourly price columns (B2:Z2), plus
# the slot summary columns (AB2, AD2, AE2, AF2), to the newly
# fetched values from the automated price update.

$newDate = Get-Date -Year 2025 -Month 8 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("A2").Value = $newDate.Date

$ws.Range("B2").Value = 109.1
$ws.Range("C2").Value = 102.84
$ws.Range("D2").Value = 101.12
$ws.Range("E2").Value = 97.54000000000001
$ws.Range("F2").Value = 93.05
$ws.Range("G2").Value = 93.05
$ws.Range("H2").Value = 102.84
$ws.Range("I2").Value = 113.12
$ws.Range("J2").Value = 113.64
$ws.Range("K2").Value = 101.63
$ws.Range("L2").Value = 85.20999999999999
$ws.Range("M2").Value = 77.75
$ws.Range("N2").Value = 63.65
$ws.Range("O2").Value = 49.85
$ws.Range("P2").Value = 50.62
$ws.Range("Q2").Value = 50.62
$ws.Range("R2").Value = 51.14
$ws.Range("S2").Value = 51.4
$ws.Range("T2").Value = 64.90000000000001
$ws.Range("U2").Value = 99
$ws.Range("V2").Value = 109.1
$ws.Range("W2").Value = 114.32
$ws.Range("X2").Value = 113.63
$ws.Range("Y2").Value = 103.57
$ws.Range("Z2").Value = 88.03

$ws.Range("AB2").Value = 110.16
$ws.Range("AD2").Value = 111.71
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 108.6
